$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Fields" row for every table) used to read "rowid" for every
# table's first field; rename that field to "id" across all tables.
$ws.Range("B2:I2").Value = "id"

# New field "bits" added to the Arch table's field list (column B, row 4).
$ws.Range("B4").Value = "bits"

# Update the active selection to match the author's cursor position.
$ws.Range("H14").Select()
